# Apply updated "dSF" (column F) values after a data repull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 1
    3  = -3
    4  = -4
    6  = -1
    8  = 1
    9  = 1
    10 = -6
    11 = -4
    12 = -1
    13 = -5
    14 = -1
    15 = -5
    19 = 4
    20 = 3
    21 = -4
    22 = 0
    23 = -6
    24 = 4
    25 = 7
    26 = -2
    27 = 1
    30 = -5
    31 = 2
    32 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
